$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: new text for column D (Price) and column E (Volume(1h)).
# $null means that column is unchanged for that row (not touched by the update).
$updates = @(
    @{ Row=2; D='68.310.94'; E='  +0.93%  ' },
    @{ Row=3; D='3.748.78'; E=$null },
    @{ Row=4; D=$null; E='  +0.00%  ' },
    @{ Row=5; D='595.20'; E='  -0.25%  ' },
    @{ Row=6; D='167.00'; E='  -1.33%  ' },
    @{ Row=7; D='3.746.86'; E='  -0.79%  ' },
    @{ Row=9; D=$null; E='  -0.89%  ' },
    @{ Row=10; D='0.159'; E='  -3.26%  ' },
    @{ Row=11; D=$null; E='  -0.25%  ' },
    @{ Row=12; D=$null; E='  -1.01%  ' },
    @{ Row=13; D=$null; E='  -6.86%  ' },
    @{ Row=14; D='36.03'; E='  -1.29%  ' },
    @{ Row=15; D='4.377.27'; E='  -0.79%  ' },
    @{ Row=16; D='3.745.83'; E='  -1.20%  ' },
    @{ Row=17; D='68.347.62'; E='  +1.04%  ' },
    @{ Row=18; D='17.89'; E='  -3.60%  ' },
    @{ Row=19; D=$null; E='  -2.33%  ' },
    @{ Row=20; D=$null; E='  -0.13%  ' },
    @{ Row=21; D=$null; E='  +1.99%  ' },
    @{ Row=22; D='467.74'; E='  -0.21%  ' },
    @{ Row=23; D=$null; E='  -2.96%  ' },
    @{ Row=24; D='84.41'; E='  +0.70%  ' },
    @{ Row=25; D=$null; E='  -2.92%  ' },
    @{ Row=26; D=$null; E='  -0.74%  ' },
    @{ Row=27; D=$null; E='  -1.51%  ' },
    @{ Row=28; D='10.10'; E='  -1.93%  ' },
    @{ Row=29; D=$null; E='  -0.08%  ' },
    @{ Row=30; D='3.895.23'; E='  -0.93%  ' },
    @{ Row=31; D=$null; E='  -4.76%  ' },
    @{ Row=32; D='7.30'; E='  -4.14%  ' },
    @{ Row=33; D='29.84'; E='  -2.25%  ' },
    @{ Row=34; D=$null; E='  -2.18%  ' },
    @{ Row=35; D='9.19'; E='  +0.58%  ' },
    @{ Row=37; D='3.704.43'; E='  -1.01%  ' },
    @{ Row=38; D=$null; E='  -2.69%  ' },
    @{ Row=39; D='3.38'; E='  -10.07%  ' },
    @{ Row=40; D=$null; E='  +0.97%  ' },
    @{ Row=41; D='0.999'; E='  -0.79%  ' },
    @{ Row=42; D=$null; E='  -0.07%  ' },
    @{ Row=43; D=$null; E='  +0.09%  ' },
    @{ Row=46; D='8.57'; E='  -1.28%  ' },
    @{ Row=47; D='43.01'; E='  +9.62%  ' },
    @{ Row=48; D=$null; E='  -1.38%  ' },
    @{ Row=49; D='45.87'; E='  +0.00%  ' },
    @{ Row=50; D='146.60'; E='  +4.06%  ' },
    @{ Row=51; D='389.49'; E='  -1.70%  ' }
)

foreach ($item in $updates) {
    $r = $item.Row

    if ($null -ne $item.D) {
        $dCell = $ws.Cells.Item($r, 4)
        # Some price strings are plain decimals (e.g. "595.20") that Excel would
        # otherwise auto-convert to a number; force text storage so the value
        # round-trips exactly like the original inline-string cell content.
        if ($item.D -match '^[+-]?\d+(\.\d+)?$') {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $item.D
    }

    if ($null -ne $item.E) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}
